$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current data rows (2-13), columns A-T, before overwriting anything.
$snapshot = @{}
for ($r = 2; $r -le 13; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le 20; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping of source row (in the original sheet) -> destination row (after edit).
$map = @{
    2  = 12
    3  = 13
    4  = 6
    5  = 7
    6  = 10
    7  = 5
    8  = 11
    9  = 2
    10 = 3
    11 = 4
    12 = 8
    13 = 9
}

foreach ($srcRow in $map.Keys) {
    $dstRow = $map[$srcRow]
    $rowVals = $snapshot[$srcRow]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($dstRow, $c).Value2 = $rowVals[$c]
    }
}
